# Auto-generated edit script: updates cryptos Price (D) and Volume(1h) (E) columns
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a plain text value without Excel auto-converting
# number-looking strings (e.g. "254.76") into numeric values. We temporarily
# mark the cell as Text, assign the value, then restore formatting by copying
# the format from a neighboring plain cell (B2) that carries no special style,
# so no stray style attributes are left on the cell.
function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range("B2").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

# --- Price column (D): values that look like plain numbers need the text-safe helper ---
Set-TextValue "D5" '254.76'
Set-TextValue "D6" '0.629'
Set-TextValue "D7" '68.33'
Set-TextValue "D9" '0.586'
Set-TextValue "D10" '37.85'
Set-TextValue "D11" '58.78'
Set-TextValue "D12" '0.0938'
Set-TextValue "D13" '7.20'
Set-TextValue "D16" '0.875'
Set-TextValue "D22" '72.11'
Set-TextValue "D23" '232.87'
Set-TextValue "D25" '3.92'
Set-TextValue "D26" '11.88'
Set-TextValue "D28" '2.53'
Set-TextValue "D30" '169.95'
Set-TextValue "D31" '20.65'
Set-TextValue "D33" '0.0746'
Set-TextValue "D34" '0.124'
Set-TextValue "D35" '5.49'
Set-TextValue "D36" '27.34'
Set-TextValue "D39" '0.0301'
Set-TextValue "D40" '12.89'
Set-TextValue "D43" '64.42'
Set-TextValue "D44" '4.95'
Set-TextValue "D45" '0.202'
Set-TextValue "D46" '8.63'
Set-TextValue "D50" '1.17'

# --- Price column (D): values that already contain multiple dots / special glyphs ---
# are never auto-converted to numbers by Excel, so a direct assignment is safe.
$ws.Range("D2").Value = '41.346.19'
$ws.Range("D3").Value = '2.188.18'
$ws.Range("D15").Value = '2.514.17'
$ws.Range("D18").Value = '2.188.05'
$ws.Range("D19").Value = '41.277.13'
$ws.Range("D20").Value = '0.0₃0956'

# --- Volume(1h) column (E): percentage strings with padding spaces are always text ---
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("E5").Value = '  +6.53%  '
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +9.38%  '
$ws.Range("E10").Value = '  +3.83%  '
$ws.Range("E11").Value = '  +2.21%  '
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("E13").Value = '  +10.47%  '
$ws.Range("E14").Value = '  +1.14%  '
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("E16").Value = '  +5.49%  '
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("E21").Value = '  +3.03%  '
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("E24").Value = '  +3.33%  '
$ws.Range("E25").Value = '  +9.24%  '
$ws.Range("E26").Value = '  +23.01%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +6.31%  '
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("E32").Value = '  +1.96%  '
$ws.Range("E33").Value = '  +7.06%  '
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("E35").Value = '  +8.80%  '
$ws.Range("E36").Value = '  +18.83%  '
$ws.Range("E37").Value = '  +10.91%  '
$ws.Range("E38").Value = '  +2.05%  '
$ws.Range("E39").Value = '  +13.56%  '
$ws.Range("E40").Value = '  +27.20%  '
$ws.Range("E41").Value = '  -1.20%  '
$ws.Range("E42").Value = '  -1.45%  '
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("E45").Value = '  +5.82%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").Value = '  +3.77%  '
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("E49").Value = '  +5.86%  '
$ws.Range("E50").Value = '  +1.77%  '
$ws.Range("E51").Value = '  -3.79%  '

$wb.Save()
